$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.279.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.664.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5325"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07847"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.551"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.668.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.893.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5524"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8185"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.677"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.054"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.010"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.39%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.245"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.472"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05854"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.582"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.311"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.618"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.822"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9592"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.431"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5813"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01616"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.885"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8537"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.009"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.047.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.806.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").Value = "  +1.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4372"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.944"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05165"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.447"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.48%  "
